$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = '10 minutes'
$ws.Range("F4").Value = 'Installed the required libraries and initiated the DB'

$ws.Range("A5").Value = 'Embeddings'
$ws.Range("B5").Value = 'Create embeddings using LLM/Sentence Transformer'
$ws.Range("E5").Value = '10 minutes'
$ws.Range("F5").Value = 'Created two approaches to get the embedding. One with Open AI embedding model and the another with Sentence Transformer. Since Open AI''s embedding model is paid, we can go with the Sentence Transformer which is free of cost.'

$ws.Range("A6").Value = 'ChatBot'
$ws.Range("B6").Value = 'Create ChatBot to give user query'

$ws.Range("A7").Value = 'Similarity check'
$ws.Range("B7").Value = ' Check the similarity of the embeddings using Cosine similarity'
$ws.Range("E7").Value = '10 minutes'
$ws.Range("F7").Value = 'ChromaDB''s default way of checking similarity is using Cosine Similarity.'

$ws.Range("A8").Value = 'Final output'
$ws.Range("B8").Value = 'Get the final output from after similarity check'
$ws.Range("E8").Value = '10 minutes'
$ws.Range("F8").Value = 'Done in Postman'

$ws.Range("A9").Value = 'Generate summary'
$ws.Range("B9").Value = 'Generate the summary of the output'

$ws.Range("A10").Value = 'Build API for upload'
$ws.Range("B10").Value = 'Endpoints for uploading and storing chunks'
$ws.Range("E10").Value = '10 minutes'
$ws.Range("F10").Value = 'For creating and testing using Postman'

$ws.Range("A11").Value = 'Build API for similarity check'
$ws.Range("B11").Value = 'Similarity search endpoint to query documents'
$ws.Range("E11").Value = '10 minutes'
$ws.Range("F11").Value = 'For creating and testing using Postman'

$ws.Range("A12").Value = 'Build API for get'
$ws.Range("B12").Value = 'endpoint to get details of specific journal when requested'
$ws.Range("E12").Value = '30 minutes'
$ws.Range("F12").Value = 'Here I was facing some issues to get the details of the specific source_doc_id requested. But was able to debug it later.'
